$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 13 for the "Docentes responsaveis" entries
# (shifts old rows 13-21 down to 15-23)
$ws.Rows("13:14").Insert()

# The inserted rows only carry column-A formatting from the row above;
# copy the B:C cell formatting used throughout the sheet onto the new cells.
$ws.Range("B10:C10").Copy()
$ws.Range("B13:C14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fix row 10 (Objetivos) - was showing the wrong duplicated text
$ws.Range("B10").Value = 'Apresentar aos alunos a Engenharia Bioquímica, as características da profissão e orientar quanto as atribuições e as áreas de atuação do Engenheiro Bioquímico. Além disso, desenvolver nos alunos uma visão macro dos tipos e etapas de um bioprocesso industrial e, por fim, orientar sobre a atuação do Engenheiro Bioquímico na indústria, pesquisa e ensino, e empreendedorismo e inovação em engenharia.'
$ws.Range("C10").Value = 'Apresentar aos alunos a Engenharia Bioquímica, as características da profissão e orientar quanto as atribuições e as áreas de atuação do Engenheiro Bioquímico. Além disso, desenvolver nos alunos uma visão macro dos tipos e etapas de um bioprocesso industrial e, por fim, orientar sobre a atuação do Engenheiro Bioquímico na indústria, pesquisa e ensino, e empreendedorismo e inovação em engenharia.'

# Fill the two new rows 13 & 14 with "Docentes responsaveis" names
$ws.Range("B13").Value = '101761 - Arnaldo Márcio Ramalho Prata'
$ws.Range("C13").Value = '101761 - Arnaldo Márcio Ramalho Prata'
$ws.Range("B14").Value = '5817181 - Valdeir Arantes'
$ws.Range("C14").Value = '5817181 - Valdeir Arantes'

# Fix row 15 (Programa resumido) - was showing stale duplicated text
$ws.Range("B15").Value = '1. Histórico da Engenharia Bioquímica 2. Engenharia Bioquímica: Definições e conceitos 3. Mercado de trabalho de Engenharia4. Áreas de atuação do Engenheiro Bioquímico5. A Indústria de Bioprocessos 6. Escalas de produção 7. Estudo de casos (processos biotecnológicos)8. Visita supervisionada.'
$ws.Range("C15").Value = '1. Histórico da Engenharia Bioquímica 2. Engenharia Bioquímica: Definições e conceitos 3. Mercado de trabalho de Engenharia4. Áreas de atuação do Engenheiro Bioquímico5. A Indústria de Bioprocessos 6. Escalas de produção 7. Estudo de casos (processos biotecnológicos)8. Visita supervisionada.'

# Fix row 17 (Programa) - was showing stale duplicated text
$ws.Range("B17").Value = '1.Histórico da Engenharia Bioquímica: interação entre ciências biológicas e a engenha, multidisciplinaridade, peculiaridades dos processos biotecnológicos. 2.Mercado de trabalho da Engenharia do Brasil 3.Atribuições e áreas de atuação do Engenheiro Bioquímico 4.Definições e conceitos – processo enzimático, processo fermentativo genérico, agentes de transformação, biorreator, matéria prima, tipos de substratos, conversão de substrato em produto, tipos de produtos biotecnológicos, recuperação de produtos, entre outros. 5.A Indústria de Bioprocessos – tipos de indústrias, equipamentos, instalações, principais operações unitárias. 6.Escalas de produção – laboratório, piloto, industrial. 7.Estudo de casos (processos biotecnológicos). 8.Empreendedorismo e Inovação em Engenharia.9.Visitas supervisionadas – visitas a laboratórios e a indústria de bioprocesso.'
$ws.Range("C17").Value = '1.Histórico da Engenharia Bioquímica: interação entre ciências biológicas e a engenha, multidisciplinaridade, peculiaridades dos processos biotecnológicos. 2.Mercado de trabalho da Engenharia do Brasil 3.Atribuições e áreas de atuação do Engenheiro Bioquímico 4.Definições e conceitos – processo enzimático, processo fermentativo genérico, agentes de transformação, biorreator, matéria prima, tipos de substratos, conversão de substrato em produto, tipos de produtos biotecnológicos, recuperação de produtos, entre outros. 5.A Indústria de Bioprocessos – tipos de indústrias, equipamentos, instalações, principais operações unitárias. 6.Escalas de produção – laboratório, piloto, industrial. 7.Estudo de casos (processos biotecnológicos). 8.Empreendedorismo e Inovação em Engenharia.9.Visitas supervisionadas – visitas a laboratórios e a indústria de bioprocesso.'

# Fix row 20 (Metodo) - was showing stale duplicated text
$ws.Range("B20").Value = 'O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras; exercícios individuais realizados no decorrer da disciplina; exercícios; dinâmicas. Para os projetos, os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a aplicações dos conceitos abordados à um processo, produto ou serviço na área de Engenharia de Bioquímica e que relacione com a formação acadêmica e atribuições profissionais do Engenheiro Bioquímico.'
$ws.Range("C20").Value = 'O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras; exercícios individuais realizados no decorrer da disciplina; exercícios; dinâmicas. Para os projetos, os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a aplicações dos conceitos abordados à um processo, produto ou serviço na área de Engenharia de Bioquímica e que relacione com a formação acadêmica e atribuições profissionais do Engenheiro Bioquímico.'

# Fix row 21 (Criterio) - was showing stale duplicated text
$ws.Range("B21").Value = 'A nota (N) será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, Avaliação dos Pares, Apresentação de Trabalhos, dentre outros.'
$ws.Range("C21").Value = 'A nota (N) será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, Avaliação dos Pares, Apresentação de Trabalhos, dentre outros.'

# Fix row 22 (Norma de recuperacao) - was showing stale duplicated text
$ws.Range("B22").Value = 'Média Final = (N + Prova Recuperação)/2'
$ws.Range("C22").Value = 'Média Final = (N + Prova Recuperação)/2'

# Fix row 23 (Bibliografia) - was showing stale duplicated text
$ws.Range("B23").Value = 'Schmidell, W.; Lima, U. A.; Aquarone, E.; Borzani, W. Biotecnologia Industrial – EngenhariaBioquímica, vol. 2, São Paulo: Edgard Blücher, 2001.Shuler, L. M.; Kargi, F. Bioprocess Engineering – Basic Concepts. Second edition. NewJersey: PrenticeHall,2002.Arigos atuais relacionaos com o tema de Engenharia Bioquímica'
$ws.Range("C23").Value = 'Schmidell, W.; Lima, U. A.; Aquarone, E.; Borzani, W. Biotecnologia Industrial – EngenhariaBioquímica, vol. 2, São Paulo: Edgard Blücher, 2001.Shuler, L. M.; Kargi, F. Bioprocess Engineering – Basic Concepts. Second edition. NewJersey: PrenticeHall,2002.Arigos atuais relacionaos com o tema de Engenharia Bioquímica'

Write-Host "Edit complete"